$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "{'metric': 'euclidean', 'n_neighbors': 39, 'weights': 'uniform'}"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "65.35%"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "66.25%"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "65.35%"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "64.85%"
